$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 is new; ensure the date cell uses the same date-time number format as D2:D6
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Terminal La Palmera de La Serena"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44454
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = "Arándano (blue)"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("Q2").Value = "$/bandeja 2 kilos"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 6750
$ws.Range("T2").Value = 2

# Row 3
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44435
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 19500
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19750
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 9875
$ws.Range("T3").Value = 2

# Row 4
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44445
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("Q4").Value = "$/bandeja 2 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 7250
$ws.Range("T4").Value = 2

# Row 5
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44446
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101001
$ws.Range("J5").Value = "Arándano (blue)"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("Q5").Value = "$/bandeja 2 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 7250
$ws.Range("T5").Value = 2

# Row 6
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44452
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13500
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 6750
$ws.Range("T6").Value = 2

# Row 7
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Terminal La Palmera de La Serena"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44448
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101001
$ws.Range("J7").Value = "Arándano (blue)"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("Q7").Value = "$/bandeja 2 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 7250
$ws.Range("T7").Value = 2
